$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '67.804.16'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = "'" + '3.811.09'
$ws.Range('E3').Value = '  +0.76%  '
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').Value = "'" + '604.09'
$ws.Range('E5').Value = '  +1.44%  '
$ws.Range('D6').Value = "'" + '166.11'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +0.34%  '
$ws.Range('E10').Value = '  +1.02%  '
$ws.Range('E11').Value = '  +0.54%  '
$ws.Range('D12').Value = "'" + '0.0000250'
$ws.Range('E12').Value = '  -0.86%  '
$ws.Range('D13').Value = "'" + '35.98'
$ws.Range('E13').Value = '  +0.12%  '
$ws.Range('D14').Value = "'" + '4.452.35'
$ws.Range('E14').Value = '  +0.76%  '
$ws.Range('D15').Value = "'" + '3.809.69'
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('D16').Value = "'" + '67.837.03'
$ws.Range('D17').Value = "'" + '18.38'
$ws.Range('E17').Value = '  -0.59%  '
$ws.Range('D18').Value = "'" + '7.09'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('E19').Value = '  +1.86%  '
$ws.Range('D20').Value = "'" + '464.10'
$ws.Range('D21').Value = "'" + '9.85'
$ws.Range('E21').Value = '  -1.36%  '
$ws.Range('D22').Value = "'" + '0.702'
$ws.Range('E22').Value = '  +0.97%  '
$ws.Range('E23').Value = '  -3.35%  '
$ws.Range('D24').Value = "'" + '83.39'
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('D25').Value = "'" + '12.16'
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('D26').Value = "'" + '2.12'
$ws.Range('E26').Value = '  -0.20%  '
$ws.Range('D27').Value = "'" + '10.03'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('D29').Value = "'" + '3.961.10'
$ws.Range('E29').Value = '  +0.79%  '
$ws.Range('E30').Value = '  +0.70%  '
$ws.Range('E31').Value = '  +2.99%  '
$ws.Range('E32').Value = '  -0.59%  '
$ws.Range('D33').Value = "'" + '29.47'
$ws.Range('E33').Value = '  -0.40%  '
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('D35').Value = "'" + '9.08'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').Value = "'" + '0.0998'
$ws.Range('E36').Value = '  -0.26%  '
$ws.Range('E37').Value = '  +0.54%  '
$ws.Range('B38').Value = 'Mantle'
$ws.Range('C38').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D38').Value = "'" + '0.997'
$ws.Range('E38').Value = '  +0.29%  '
$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').Value = "'" + '5.82'
$ws.Range('E39').Value = '  +1.26%  '
$ws.Range('D40').Value = "'" + '3.23'
$ws.Range('E40').Value = '  -3.04%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E43').Value = '  -2.74%  '
$ws.Range('D44').Value = "'" + '47.73'
$ws.Range('E44').Value = '  -0.93%  '
$ws.Range('E45').Value = '  +0.44%  '
$ws.Range('D46').Value = "'" + '28.20'
$ws.Range('E46').Value = '  +6.54%  '
$ws.Range('D47').Value = "'" + '151.68'
$ws.Range('E47').Value = '  +1.51%  '
$ws.Range('E48').Value = '  +12.22%  '
$ws.Range('E49').Value = '  +0.51%  '
$ws.Range('E50').Value = '  +1.67%  '
$ws.Range('D51').Value = "'" + '390.80'
$ws.Range('E51').Value = '  -0.35%  '
